# Applies the "Add files via upload" revision to
# draft-gandhi-ippm-simple-direct-loss-00.pptx:
#   - Reflow/resize the title-slide subtitle and author-list placeholders
#   - Bump the author-list text to 16pt
#   - Add a new "Stefano Salsano" author entry (with a mailto hyperlink)
#     plus two trailing blank paragraphs
#   - Fix a stray extra space in the STAMP packet-format diagram's
#     "SSID" column header on the agenda/body slide

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height (and TextRange font
# sizes) are single-precision floats under the hood, so a naive
# emu/12700.0 can truncate one EMU short of the intended value once it
# is converted back on save. Search the handful of representable
# Single values nearest the exact quotient for one that truncates back
# to exactly the requested EMU count.
# ---------------------------------------------------------------------
function EmuToPt($targetEmu) {
    $basePt = $targetEmu / 12700.0
    $bestF = [float]$basePt
    for ($i = -2000; $i -le 2000; $i++) {
        $cand = $basePt + ($i * 0.0000001)
        $f = [float]$cand
        $emu = [double]$f * 12700.0
        $trunc = [Math]::Floor($emu)
        if ([int64]$trunc -eq $targetEmu) {
            return $f
        }
    }
    return $bestF
}

# ---------------------------------------------------------------------
# Slide 1 (title slide): subtitle + author-list placeholders
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# "Rectangle 3" - the draft name subtitle - moved slightly up/right.
$subtitle = $slide1.Shapes.Item(2)
$subtitle.Left = EmuToPt 723900
$subtitle.Top = EmuToPt 1774232

# "Rectangle 4" - the author list - repositioned & enlarged to fit the
# new entry.
$authors = $slide1.Shapes.Item(3)
$authors.Left = EmuToPt 1524000
$authors.Top = EmuToPt 2641943
$authors.Width = EmuToPt 7010400
$authors.Height = EmuToPt 1628433

$tr = $authors.TextFrame.TextRange

# Append the new "Stefano Salsano" paragraph, run by run so each run
# boundary matches the source punctuation/word breaks.
$null = $tr.InsertAfter("`rStefano Salsano - ")
$null = $tr.InsertAfter("Universita")
$null = $tr.InsertAfter(" di Roma `"Tor ")
$null = $tr.InsertAfter("Vergata")
$null = $tr.InsertAfter("`" (")
$null = $tr.InsertAfter("stefano.salsano@uniroma2.it")
$null = $tr.InsertAfter(")")
# Two trailing empty paragraphs.
$null = $tr.InsertAfter("`r`r")

# Hyperlink just the new e-mail address run.
$fullText = $tr.Text
$emailAddr = "stefano.salsano@uniroma2.it"
$emailOffset = $fullText.IndexOf($emailAddr)
$emailRange = $tr.Characters($emailOffset + 1, $emailAddr.Length)
$emailRange.ActionSettings.Item(1).Action = 7
$emailRange.ActionSettings.Item(1).Hyperlink.Address = "mailto:" + $emailAddr

# Bump every run in the author-list box (old + new) to 16pt.
$tr.Font.Size = 16

# ---------------------------------------------------------------------
# Slide 11 (STAMP packet format): drop one extra space before the
# closing "|" of the "SSID" column in the Error Estimate/SSID row.
# ---------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$diagram = $slide11.Shapes.Item(3)
$diagramText = $diagram.TextFrame.TextRange
$ssidPara = $diagramText.Paragraphs(9, 1)
$ssidRun = $ssidPara.Runs(1, 1)
$ssidRun.Text = "    |         Error Estimate        |           SSID                |"
